$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every value (including dates like "2024-05-14" and
# times like "11:20:09") as literal text, not as Excel date/time serial
# numbers. Pre-format the new range as Text so values like "2024-05-14"
# are not auto-converted to dates when assigned via .Value.
$newRange = $ws.Range("A130:G158")
$newRange.NumberFormat = "@"

# Row 130
$ws.Range("A130").Value = '2024-05-14'
$ws.Range("B130").Value = '11:20:09'
$ws.Range("C130").Value = 'Palet atascado en la curva'
$ws.Range("D130").Value = '-'
$ws.Range("E130").Value = '-'
$ws.Range("F130").Value = '-'
$ws.Range("G130").Value = '-'

# Row 131
$ws.Range("A131").Value = '2024-05-14'
$ws.Range("B131").Value = '11:20:16'
$ws.Range("C131").Value = 'AOI no detecta pieza'
$ws.Range("D131").Value = '-'
$ws.Range("E131").Value = '-'
$ws.Range("F131").Value = '-'
$ws.Range("G131").Value = '-'

# Row 132
$ws.Range("A132").Value = '2024-05-14'
$ws.Range("B132").Value = '11:20:27'
$ws.Range("C132").Value = '-'
$ws.Range("D132").Value = 'Cámara no detecta skeleton'
$ws.Range("E132").Value = '-'
$ws.Range("F132").Value = '-'
$ws.Range("G132").Value = '-'

# Row 133
$ws.Range("A133").Value = '2024-05-14'
$ws.Range("B133").Value = '11:22:33'
$ws.Range("C133").Value = 'Secuencia atornillador'
$ws.Range("D133").Value = '-'
$ws.Range("E133").Value = '-'
$ws.Range("F133").Value = '-'
$ws.Range("G133").Value = '-'

# Row 134
$ws.Range("A134").Value = '2024-05-14'
$ws.Range("B134").Value = '11:26:32'
$ws.Range("C134").Value = 'Palet atascado en la curva'
$ws.Range("D134").Value = '-'
$ws.Range("E134").Value = '-'
$ws.Range("F134").Value = '-'
$ws.Range("G134").Value = '-'

# Row 135
$ws.Range("A135").Value = '2024-05-14'
$ws.Range("B135").Value = '11:26:50'
$ws.Range("C135").Value = '-'
$ws.Range("D135").Value = '-'
$ws.Range("E135").Value = '-'
$ws.Range("F135").Value = 'Pieza enganchada en HV Test'
$ws.Range("G135").Value = '-'

# Row 136
$ws.Range("A136").Value = '2024-05-14'
$ws.Range("B136").Value = '11:26:53'
$ws.Range("C136").Value = '-'
$ws.Range("D136").Value = '-'
$ws.Range("E136").Value = '-'
$ws.Range("F136").Value = 'Core enganchado'
$ws.Range("G136").Value = '-'

# Row 137
$ws.Range("A137").Value = '2024-05-14'
$ws.Range("B137").Value = '11:26:56'
$ws.Range("C137").Value = '-'
$ws.Range("D137").Value = '-'
$ws.Range("E137").Value = '-'
$ws.Range("F137").Value = 'Traza'
$ws.Range("G137").Value = '-'

# Row 138
$ws.Range("A138").Value = '2024-05-14'
$ws.Range("B138").Value = '11:26:59'
$ws.Range("C138").Value = '-'
$ws.Range("D138").Value = '-'
$ws.Range("E138").Value = '-'
$ws.Range("F138").Value = 'Fallo cámara QR'
$ws.Range("G138").Value = '-'

# Row 139
$ws.Range("A139").Value = '2024-05-14'
$ws.Range("B139").Value = '11:30:59'
$ws.Range("C139").Value = 'No pone tornillo'
$ws.Range("D139").Value = '-'
$ws.Range("E139").Value = '-'
$ws.Range("F139").Value = '-'
$ws.Range("G139").Value = '-'

# Row 140
$ws.Range("A140").Value = '2024-05-14'
$ws.Range("B140").Value = '11:35:13'
$ws.Range("C140").Value = '-'
$ws.Range("D140").Value = 'Tornillo atascado en tolva'
$ws.Range("E140").Value = '-'
$ws.Range("F140").Value = '-'
$ws.Range("G140").Value = '-'

# Row 141
$ws.Range("A141").Value = '2024-05-14'
$ws.Range("B141").Value = '11:35:17'
$ws.Range("C141").Value = '-'
$ws.Range("D141").Value = '-'
$ws.Range("E141").Value = 'No lee QR'
$ws.Range("F141").Value = '-'
$ws.Range("G141").Value = '-'

# Row 142
$ws.Range("A142").Value = '2024-05-14'
$ws.Range("B142").Value = '11:35:32'
$ws.Range("C142").Value = '-'
$ws.Range("D142").Value = '-'
$ws.Range("E142").Value = 'Etiquetadora'
$ws.Range("F142").Value = '-'
$ws.Range("G142").Value = '-'

# Row 143
$ws.Range("A143").Value = '2024-05-14'
$ws.Range("B143").Value = '11:35:37'
$ws.Range("C143").Value = '-'
$ws.Range("D143").Value = '-'
$ws.Range("E143").Value = '-'
$ws.Range("F143").Value = 'Robot no coloca bien filter en palet'
$ws.Range("G143").Value = '-'

# Row 144
$ws.Range("A144").Value = '2024-05-14'
$ws.Range("B144").Value = '11:35:43'
$ws.Range("C144").Value = '-'
$ws.Range("D144").Value = '-'
$ws.Range("E144").Value = '-'
$ws.Range("F144").Value = 'Core enganchado'
$ws.Range("G144").Value = '-'

# Row 145
$ws.Range("A145").Value = '2024-05-14'
$ws.Range("B145").Value = '11:39:49'
$ws.Range("C145").Value = 'Ascensor no sube'
$ws.Range("D145").Value = '-'
$ws.Range("E145").Value = '-'
$ws.Range("F145").Value = '-'
$ws.Range("G145").Value = '-'

# Row 146
$ws.Range("A146").Value = '2024-05-14'
$ws.Range("B146").Value = '11:42:40'
$ws.Range("C146").Value = 'Secuencia atornillador'
$ws.Range("D146").Value = '-'
$ws.Range("E146").Value = '-'
$ws.Range("F146").Value = '-'
$ws.Range("G146").Value = '-'

# Row 147
$ws.Range("A147").Value = '2024-05-14'
$ws.Range("B147").Value = '11:43:49'
$ws.Range("C147").Value = '-'
$ws.Range("D147").Value = '-'
$ws.Range("E147").Value = '-'
$ws.Range("F147").Value = 'Robot no coge PCB'
$ws.Range("G147").Value = '-'

# Row 148
$ws.Range("A148").Value = '2024-05-14'
$ws.Range("B148").Value = '11:43:55'
$ws.Range("C148").Value = '-'
$ws.Range("D148").Value = '-'
$ws.Range("E148").Value = '-'
$ws.Range("F148").Value = 'Cover atascado'
$ws.Range("G148").Value = '-'

# Row 149
$ws.Range("A149").Value = '2024-05-14'
$ws.Range("B149").Value = '11:45:47'
$ws.Range("C149").Value = '-'
$ws.Range("D149").Value = '-'
$ws.Range("E149").Value = 'Power atascado en prensa, cuesta sacar'
$ws.Range("F149").Value = '-'
$ws.Range("G149").Value = '-'

# Row 150
$ws.Range("A150").Value = '2024-05-14'
$ws.Range("B150").Value = '11:47:03'
$ws.Range("C150").Value = 'Fallo en paletizador'
$ws.Range("D150").Value = '-'
$ws.Range("E150").Value = '-'
$ws.Range("F150").Value = '-'
$ws.Range("G150").Value = '-'

# Row 151
$ws.Range("A151").Value = '2024-05-14'
$ws.Range("B151").Value = '11:58:11'
$ws.Range("C151").Value = '-'
$ws.Range("D151").Value = '-'
$ws.Range("E151").Value = 'La cámara no detecta Busbar'
$ws.Range("F151").Value = '-'
$ws.Range("G151").Value = '-'

# Row 152
$ws.Range("A152").Value = '2024-05-14'
$ws.Range("B152").Value = '12:09:25'
$ws.Range("C152").Value = 'Fallo en paletizador'
$ws.Range("D152").Value = '-'
$ws.Range("E152").Value = '-'
$ws.Range("F152").Value = '-'
$ws.Range("G152").Value = '-'

# Row 153
$ws.Range("A153").Value = '2024-05-14'
$ws.Range("B153").Value = '12:16:33'
$ws.Range("C153").Value = '-'
$ws.Range("D153").Value = 'AOI (malla)'
$ws.Range("E153").Value = '-'
$ws.Range("F153").Value = '-'
$ws.Range("G153").Value = '-'

# Row 154
$ws.Range("A154").Value = '2024-05-14'
$ws.Range("B154").Value = '12:16:37'
$ws.Range("C154").Value = '-'
$ws.Range("D154").Value = 'Detección de sealling mal puesto'
$ws.Range("E154").Value = '-'
$ws.Range("F154").Value = '-'
$ws.Range("G154").Value = '-'

# Row 155
$ws.Range("A155").Value = '2024-05-14'
$ws.Range("B155").Value = '12:16:48'
$ws.Range("C155").Value = '-'
$ws.Range("D155").Value = '-'
$ws.Range("E155").Value = '-'
$ws.Range("F155").Value = 'Repeat funcional'
$ws.Range("G155").Value = '-'

# Row 156
$ws.Range("A156").Value = '2024-05-14'
$ws.Range("B156").Value = '12:16:52'
$ws.Range("C156").Value = '-'
$ws.Range("D156").Value = '-'
$ws.Range("E156").Value = '-'
$ws.Range("F156").Value = 'Core enganchado'
$ws.Range("G156").Value = '-'

# Row 157
$ws.Range("A157").Value = '2024-05-14'
$ws.Range("B157").Value = '12:16:57'
$ws.Range("C157").Value = '-'
$ws.Range("D157").Value = '-'
$ws.Range("E157").Value = '-'
$ws.Range("F157").Value = 'Robot no coloca bien filter en palet'
$ws.Range("G157").Value = '-'

# Row 158
$ws.Range("A158").Value = '2024-05-14'
$ws.Range("B158").Value = '12:17:01'
$ws.Range("C158").Value = '-'
$ws.Range("D158").Value = '-'
$ws.Range("E158").Value = '-'
$ws.Range("F158").Value = 'Traza'
$ws.Range("G158").Value = '-'

# Drop the temporary Text number format so the new cells end up with
# the same default (unstyled) formatting as the rest of the sheet,
# while keeping the literal text values just assigned above.
$newRange.ClearFormats()
